# Refresh cryptocurrency price/volume snapshot (Price column D, Volume(1h) column E).
# Values are written as plain text to match the inline-string cells already used
# in the sheet; a leading apostrophe is used for numeric-looking Price values so
# Excel keeps them as text instead of auto-converting to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.885.98'
$ws.Range('E2').Value = '  -1.78%  '
$ws.Range('D3').Value = '3.515.25'
$ws.Range('E3').Value = '  -3.56%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '''579.88'
$ws.Range('E5').Value = '  -4.36%  '
$ws.Range('E6').Value = '  -4.75%  '
$ws.Range('D7').Value = '3.503.59'
$ws.Range('E7').Value = '  -3.49%  '
$ws.Range('E8').Value = '  -2.54%  '
$ws.Range('D9').Value = '''1.00'
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').Value = '''0.206'
$ws.Range('E10').Value = '  -6.82%  '
$ws.Range('E11').Value = '  -4.27%  '
$ws.Range('D12').Value = '''51.52'
$ws.Range('E12').Value = '  -4.49%  '
$ws.Range('E13').Value = '  -6.57%  '
$ws.Range('E14').Value = '  -4.57%  '
$ws.Range('D15').Value = '4.069.24'
$ws.Range('E15').Value = '  -3.63%  '
$ws.Range('D16').Value = '''648.46'
$ws.Range('E16').Value = '  -4.39%  '
$ws.Range('D17').Value = '69.854.04'
$ws.Range('E17').Value = '  -1.78%  '
$ws.Range('D18').Value = '3.522.00'
$ws.Range('E18').Value = '  -3.35%  '
$ws.Range('D19').Value = '''12.45'
$ws.Range('E19').Value = '  -3.71%  '
$ws.Range('D20').Value = '''18.41'
$ws.Range('E20').Value = '  -3.67%  '
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('E22').Value = '  -4.99%  '
$ws.Range('D23').Value = '''18.13'
$ws.Range('E23').Value = '  -4.46%  '
$ws.Range('E24').Value = '  -0.84%  '
$ws.Range('D25').Value = '''99.68'
$ws.Range('E25').Value = '  -5.36%  '
$ws.Range('D26').Value = '''4.31'
$ws.Range('E26').Value = '  -7.04%  '
$ws.Range('E27').Value = '  -4.48%  '
$ws.Range('D28').Value = '''10.15'
$ws.Range('E28').Value = '  -4.55%  '
$ws.Range('D29').Value = '''9.43'
$ws.Range('E29').Value = '  -5.29%  '
$ws.Range('D30').Value = '''32.83'
$ws.Range('E30').Value = '  -4.74%  '
$ws.Range('D31').Value = '''4.29'
$ws.Range('E31').Value = '  -6.11%  '
$ws.Range('E32').Value = '  -6.90%  '
$ws.Range('D33').Value = '''11.67'
$ws.Range('E33').Value = '  -4.87%  '
$ws.Range('D34').Value = '''574.73'
$ws.Range('E34').Value = '  +10.21%  '
$ws.Range('E35').Value = '  -4.72%  '
$ws.Range('D36').Value = '''61.45'
$ws.Range('E36').Value = '  -3.17%  '
$ws.Range('D37').Value = '3.780.41'
$ws.Range('E38').Value = '  +0.08%  '
$ws.Range('D39').Value = '0.0₃0796'
$ws.Range('E39').Value = '  -9.32%  '
$ws.Range('E40').Value = '  -0.27%  '
$ws.Range('D41').Value = '''3.76'
$ws.Range('E41').Value = '  +41.27%  '
$ws.Range('E42').Value = '  -4.67%  '
$ws.Range('D43').Value = '''0.376'
$ws.Range('E43').Value = '  -4.15%  '
$ws.Range('E44').Value = '  -2.45%  '
$ws.Range('D45').Value = '''34.36'
$ws.Range('E45').Value = '  -6.48%  '
$ws.Range('E46').Value = '  -3.08%  '
$ws.Range('D47').Value = '''2.89'
$ws.Range('E47').Value = '  -5.73%  '
$ws.Range('E48').Value = '  -3.24%  '
$ws.Range('E49').Value = '  -3.52%  '
$ws.Range('E50').Value = '  -0.22%  '
$ws.Range('D51').Value = '''8.23'
$ws.Range('E51').Value = '  -5.21%  '
